$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Tree Name"
$ws.Range("B1").Value = "Species"
$ws.Range("C1").Value = "CO2_per_year_kg"
$ws.Range("D1").Value = "Max_Age"
$ws.Range("E1").Value = "Common_Location"

# Data rows
$data = @(
    @("Neem",       "Azadirachta indica",       26, 100, "East Godavari"),
    @("Indian Almond","Terminalia catappa",     30, 80,  "Kakinada"),
    @("Pongamia",   "Millettia pinnata",        35, 60,  "Godavari Belt"),
    @("Banyan",     "Ficus benghalensis",       55, 200, "Rural & Urban Areas"),
    @("Peepal",     "Ficus religiosa",          48, 150, "Village Roads"),
    @("Teak",       "Tectona grandis",          22, 80,  "Plantations & Farms"),
    @("Jamun",      "Syzygium cumini",          28, 90,  "Home Gardens"),
    @("Casuarina",  "Casuarina equisetifolia",  20, 40,  "Coastal Belt"),
    @("Tamarind",   "Tamarindus indica",        32, 120, "Village Borders"),
    @("Amla",       "Phyllanthus emblica",      25, 60,  "Backyards & Schools")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}
